# Automatische test-sync: 2025-06-22 18:43:50
$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: append new row 13 ---
$logs.Range("A13").Value = "Beschadigd product ontvangen"
$logs.Range("B13").Value = "mailmind.test@zohomail.eu"
$logs.Range("C13").Value = "Het product dat ik heb ontvangen is beschadigd aangekomen."
$logs.Range("D13").Value = "Retour / Terugbetaling"
$logs.Range("E13").Value = "Beste klant,`nBedankt voor uw bericht. Wat vervelend om te horen dat het product beschadigd is aangekomen. Kunt u ons meer details geven over de schade? Bijvoorbeeld, wat voor product is het en wat voor schade is er precies aan? Eventuele foto's van de schade kunnen ook nuttig zijn.`nZodra we deze informatie hebben, zullen we ons best doen om een passende oplossing voor u te vinden.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$logs.Range("F13").Value = "2025-06-22 18:43:12"
$logs.Range("G13").Value = "Ja"

# --- Dashboard sheet: append new row 10 ---
$dash.Range("A10").Value = "Retour / Terugbetaling"
$dash.Range("B10").Value = 1

# --- Extend conditional formatting ranges to cover the new row ---
$dFcs = $logs.Range("D2:D12").FormatConditions
$dFcs.Item(1).ModifyAppliesToRange($logs.Range("D2:D13"))

$gFcs = $logs.Range("G2:G12").FormatConditions
$gFcs.Item(1).ModifyAppliesToRange($logs.Range("G2:G13"))

# --- Extend chart series ranges to include the new Dashboard row ---
$chartObj = $dash.ChartObjects(1)
$series = $chartObj.Chart.SeriesCollection(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$10,Dashboard!`$B`$2:`$B`$10,1)"
